$d = $word.ActiveDocument

# --------------------------------------------------------------------
# Helper: split the run(s) so that a boundary exists right after the
# text identified by $findText (first match only). Achieved by
# toggling a formatting property on/off on the located range - this
# forces the host to re-segment the paragraph's runs at that boundary
# without touching the text itself or merging any other split already
# made.
# --------------------------------------------------------------------
function Split-After {
    param([string]$findText)
    $r = $d.Content
    $ok = $r.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false)
    if (-not $ok) {
        Write-Output "WARNING: Split-After could not find [$findText]"
        return
    }
    $r.Bold = $true
    $r.Bold = $false
}

# ------------------------------------------------------------------
# 1. Remove the original "_GoBack" bookmark from the abstract
#    (it currently sits between "independently " and "double-screened ").
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. Abstract: "human screener" -> "human screeners"
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "conflict rate between human screener, we find",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "conflict rate between human screeners, we find", 2) | Out-Null

# ------------------------------------------------------------------
# 3. Abstract: "...) than humans. Finally, to implement this procedure we
#    have developed..." -> "...) than humans. We also have developed..."
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "than humans. Finally, to implement this procedure we",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "than humans. We also", 2) | Out-Null

# ------------------------------------------------------------------
# 4. Introduction: "independent double screening" -> "independent human
#    double screening", and the "_GoBack" bookmark now lives right after
#    the newly-inserted "human ".
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "this involves independent double screening of all references",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "this involves independent human double screening of all references", 2) | Out-Null

# ------------------------------------------------------------------
# 5. Results: drop the trailing ", and one BMJ review.    " from the
#    Campbell Systematic Reviews sentence.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "two reviews from Review of Educational Research, and one BMJ review.    ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "two reviews from Review of Educational Research", 2) | Out-Null

# ------------------------------------------------------------------
# 6. Re-create the run boundaries the original document had in those
#    paragraphs (the text edits above coalesce runs that the diff
#    shows as staying distinct). Work right-to-left within each
#    paragraph so earlier splits are not swallowed by later ones.
# ------------------------------------------------------------------

# --- Abstract paragraph ---
Split-After ") than humans. We also "
Split-After ") than humans. We "
Split-After ") than humans. W"
Split-After "highly reliable second screener, with "
Split-After "fewer false excl"
Split-After "highly reliable second"
Split-After "human screeners"
Split-After "human screener"

# --- Introduction paragraph ---
Split-After "this involves independent human "

# --- Campbell Systematic Reviews paragraph ---
Split-After "two reviews from "

# ------------------------------------------------------------------
# 7. Insert the "_GoBack" bookmark right after "human " in the
#    introduction paragraph (before "double screening").
# ------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute(
    "this involves independent human ",
    $true, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$r.Collapse(0)
$d.Bookmarks.Add("_GoBack", $r) | Out-Null
